# Revert the "Net Income" column (X) that had been added to the header row.
# This removes the shared string, the extra bold/orange header style that was
# only used by that cell, and shifts the remaining headers (EPF Employee,
# ESIC Employee, Professional Tax, Labour Welfare Fund) one column to the
# left - exactly undoing the prior merge commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column X entirely ("Net Income") - cells to its right shift left,
# the shared string is dropped when no longer referenced, dimension/row
# spans/col widths are recalculated automatically.
$ws.Columns("X:X").Delete()

# Restore the original active selection on the sheet.
$ws.Range("Z15").Select()
